$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dates = @(
    "Thu Feb 16 12:32:26 EST 2023",
    "Thu Feb 16 12:33:08 EST 2023",
    "Thu Feb 16 12:33:50 EST 2023",
    "Thu Feb 16 12:34:32 EST 2023",
    "Thu Feb 16 12:35:18 EST 2023",
    "Thu Feb 16 12:36:02 EST 2023",
    "Thu Feb 16 12:36:45 EST 2023",
    "Thu Feb 16 12:37:26 EST 2023",
    "Thu Feb 16 12:38:07 EST 2023",
    "Thu Feb 16 12:38:49 EST 2023",
    "Thu Feb 16 12:39:30 EST 2023",
    "Thu Feb 16 12:40:12 EST 2023",
    "Thu Feb 16 12:40:53 EST 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}

$wb.Save()
